$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template gets five more process blocks ("process _6" .. "process _10"),
# each laid out exactly like the existing "process _1" .. "process _5" blocks:
# a process-name header cell followed by loading time / process time /
# unloading time / cycle time sub-headers.
$headers = @(
    @("AE1", "process _6"),  @("AF1", "loading time"), @("AG1", "process time"), @("AH1", "unloading time"), @("AI1", "Cycle time"),
    @("AJ1", "process _7"),  @("AK1", "loading time"), @("AL1", "process time"), @("AM1", "unloading time"), @("AN1", "Cycle time"),
    @("AO1", "process _8"),  @("AP1", "loading time"), @("AQ1", "process time"), @("AR1", "unloading time"), @("AS1", "cycle time"),
    @("AT1", "process _9"),  @("AU1", "loading time"), @("AV1", "process time"), @("AW1", "unloading time"), @("AX1", "cycle time"),
    @("AY1", "process _10"), @("AZ1", "loading time"), @("BA1", "process time"), @("BB1", "unloading time"), @("BC1", "cycle time")
)

# Copy the existing header style (dark fill, bold white Arial, centered +
# wrapped) from A1 so the new header cells reuse the same cell style instead
# of creating new ones.
$ws.Range("A1").Copy()

foreach ($pair in $headers) {
    $addr = $pair[0]
    $text = $pair[1]
    $ws.Range($addr).Value = $text
    $ws.Range($addr).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
